$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - new row index numbers
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 6

# Column D - addresses
$ws.Range("D1").Value = "30 Wallace Green Way, Walkern, Stevenage (SG2 7FB)"
$ws.Range("D2").Value = "Vineyard Barn, The Vineyard, Welwyn Garden City (AL8 7PU)"
$ws.Range("D3").Value = "31 Wallace Green Way, Walkern, Stevenage (SG2 7FB)"
$ws.Range("D4").Value = "33 Wallace Green Way, Walkern, Stevenage (SG2 7FB)"
$ws.Range("D5").Value = "34 Wallace Green Way, Walkern, Stevenage (SG2 7FB)"
$ws.Range("D6").Value = "31 Wallace Green Way, Walkern, Stevenage (SG2 7FB)"

# Column E - numbers
$ws.Range("E1").Value = 33333
$ws.Range("E2").Value = 55556
$ws.Range("E3").Value = 55555
$ws.Range("E4").Value = 55557
$ws.Range("E5").Value = 55558
$ws.Range("E6").Value = 55555

# Column I - title numbers
$ws.Range("I1").Value = "Hg545432"
$ws.Range("I2").Value = "HD602385"
$ws.Range("I3").Value = "HD602384"
$ws.Range("I4").Value = "HD602386"
$ws.Range("I5").Value = "HD602387"
$ws.Range("I6").Value = "HD602384"

# Columns J and K are no longer used - clear them
$ws.Range("J1:J6").ClearContents()
$ws.Range("K1:K6").ClearContents()

# Column M - status
$ws.Range("M1").Value = "Submitted"
$ws.Range("M2").Value = "Submitted"
$ws.Range("M3").ClearContents()
$ws.Range("M4").Value = "Submitted"
$ws.Range("M5").Value = "Submitted"
$ws.Range("M6").ClearContents()

# Update the view: reposition selection
$ws.Range("N12").Select()
